$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing column D explanation text before overwriting, since
# the update shifts each value down one row (and the original D4 text ends
# up in the newly populated D7).
$origD4 = $ws.Range("D4").Text
$origD5 = $ws.Range("D5").Text
$origD6 = $ws.Range("D6").Text

$ws.Range("D4").Value = $origD5
$ws.Range("D5").Value = $origD6
$ws.Range("D6").Value = "More insights in the manner of the above"
$ws.Range("D7").Value = $origD4

# Update the active selection to reflect the new cell of interest.
$ws.Range("D10").Select()
